# Filter - Study - Test Suit
# Update the "startup" sheet: rename the CasesTab row label to ParticipantsTab
# and move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# A2 previously held "CasesTab" - rename it to "ParticipantsTab".
$ws.Range("A2").Value = "ParticipantsTab"

# Update the selected/active cell shown when the sheet is opened.
$ws.Activate()
$ws.Range("A2").Select()
